# Updates cryptos list prices (D) and 1h volume % changes (E)
# for the rows whose source data changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell that carries the workbook's default (unstyled) format,
# used to restore a cell's style after a NumberFormat round-trip.
$defaultStyle = $ws.Range("D4").Style

$ws.Range("D2").Value = "26.294.14"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "1.608.39"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  +0.05%  "
# D5: numeric-looking text -> force Text format so it stays a string
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.94"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("E8").Value = "  +0.67%  "
# D9: numeric-looking text -> force Text format so it stays a string
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0615"
$ws.Range("D9").Style = $defaultStyle
$ws.Range("E9").Value = "  -0.05%  "
# D10: numeric-looking text -> force Text format so it stays a string
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.44"
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").Value = "  +1.93%  "
# D11: numeric-looking text -> force Text format so it stays a string
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0813"
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").Value = "1.598.57"
$ws.Range("E13").Value = "  -0.78%  "
$ws.Range("E14").Value = "  +0.33%  "
$ws.Range("E15").Value = "  +0.89%  "
$ws.Range("D16").Value = "26.264.00"
$ws.Range("E16").Value = "  +0.49%  "
# D17: numeric-looking text -> force Text format so it stays a string
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.18"
$ws.Range("D17").Style = $defaultStyle
$ws.Range("E17").Value = "  +2.73%  "
$ws.Range("E18").Value = "  +0.77%  "
$ws.Range("E19").Value = "  +0.01%  "
# D20: numeric-looking text -> force Text format so it stays a string
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "201.72"
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = "  -0.81%  "
$ws.Range("E21").Value = "  +0.73%  "
# D22: numeric-looking text -> force Text format so it stays a string
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.33"
$ws.Range("D22").Style = $defaultStyle
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("E24").Value = "  +0.86%  "
$ws.Range("E25").Value = "  +1.19%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  -1.30%  "
$ws.Range("E28").Value = "  +0.49%  "
$ws.Range("E29").Value = "  +2.20%  "
# D30: numeric-looking text -> force Text format so it stays a string
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0498"
$ws.Range("D30").Style = $defaultStyle
$ws.Range("E30").Value = "  +5.67%  "
# D32: numeric-looking text -> force Text format so it stays a string
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.19"
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").Value = "  +2.68%  "
# D33: numeric-looking text -> force Text format so it stays a string
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.94"
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E33").Value = "  -1.50%  "
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("E35").Value = "  +1.50%  "
$ws.Range("D36").Value = "1.162.13"
$ws.Range("E36").Value = "  +3.47%  "
$ws.Range("E37").Value = "  +2.35%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("E39").Value = "  +1.03%  "
# D40: numeric-looking text -> force Text format so it stays a string
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.788"
$ws.Range("D40").Style = $defaultStyle
# D41: numeric-looking text -> force Text format so it stays a string
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.496"
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Value = "  +0.83%  "
# D42: numeric-looking text -> force Text format so it stays a string
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.36"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Value = "  +3.93%  "
# D43: numeric-looking text -> force Text format so it stays a string
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.784"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").Value = "1.743.80"
$ws.Range("E44").Value = "  +0.24%  "
# D45: numeric-looking text -> force Text format so it stays a string
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.98"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = "  -0.86%  "
$ws.Range("E46").Value = "  +14.32%  "
$ws.Range("E47").Value = "  +0.85%  "
# D48: numeric-looking text -> force Text format so it stays a string
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.06"
$ws.Range("D48").Style = $defaultStyle
$ws.Range("E48").Value = "  +1.07%  "
$ws.Range("E49").Value = "  +0.39%  "
$ws.Range("E50").Value = "  -0.30%  "
$ws.Range("E51").Value = "  -0.09%  "
